$wb = $excel.ActiveWorkbook

$alc = $wb.Worksheets.Item("ALC")
# Row 40 (ALC)
$alc.Range("H40").Value = 2410.4583
$alc.Range("I40").Value = 1494.4445
$alc.Range("J40").Value = 2960.0667
$alc.Range("K40").Value = 1494.4445
$alc.Range("L40").Value = 2960.0667
$alc.Range("M40").Value = -1319.4445
$alc.Range("N40").Value = -3310.0667

# Row 62 (ALC)
$alc.Range("H62").Value = 3167.6667
$alc.Range("I62").Value = 3000
$alc.Range("K62").Value = 3000
$alc.Range("M62").Value = -2376
$alc.Range("N62").ClearContents()

# Row 65 (ALC)
$alc.Range("H65").Value = 3167.6667
$alc.Range("I65").Value = 3000
$alc.Range("K65").Value = 15000
$alc.Range("M65").Value = -11880
$alc.Range("N65").ClearContents()

# Row 76 (ALC)
$alc.Range("H76").Value = 5625.875
$alc.Range("I76").Value = 5001.5
$alc.Range("J76").Value = 5834
$alc.Range("K76").Value = 5001.5
$alc.Range("L76").Value = 5834
$alc.Range("M76").Value = -4686.5
$alc.Range("N76").Value = -6464

# Row 79 (ALC)
$alc.Range("H79").Value = 5625.875
$alc.Range("I79").Value = 5001.5
$alc.Range("J79").Value = 5834
$alc.Range("K79").Value = 5001.5
$alc.Range("L79").Value = 5834
$alc.Range("M79").Value = -3909.5
$alc.Range("N79").Value = -8018

# Row 97 (ALC)
$alc.Range("H97").Value = 1213.75
$alc.Range("J97").Value = 1315.7142
$alc.Range("L97").Value = 3947.1426
$alc.Range("N97").Value = -4939.142599999999

# Row 112 (ALC)
$alc.Range("H112").Value = 1481.2174
$alc.Range("J112").Value = 1593.7142
$alc.Range("L112").Value = 4781.142599999999
$alc.Range("N112").Value = -6997.142599999999

# Row 129 (ALC)
$alc.Range("H129").Value = 4595.2593
$alc.Range("J129").Value = 995.85
$alc.Range("L129").Value = 2987.55
$alc.Range("N129").Value = -12987.55

# Row 137 (ALC)
$alc.Range("H137").Value = 1391.3334
$alc.Range("I137").Value = 1342.6666
$alc.Range("J137").Value = 1683.3334
$alc.Range("K137").Value = 4027.9998
$alc.Range("L137").Value = 5050.0002
$alc.Range("M137").Value = -1477.9998
$alc.Range("N137").Value = -10150.0002

# Row 138 (ALC)
$alc.Range("H138").Value = 3728.9092
$alc.Range("I138").Value = 3904.5557
$alc.Range("J138").Value = 3701.1755
$alc.Range("K138").Value = 11713.6671
$alc.Range("L138").Value = 11103.5265
$alc.Range("M138").Value = -6573.667099999999
$alc.Range("N138").Value = -21383.5265

# Row 141 (ALC)
$alc.Range("H141").Value = 5233.125
$alc.Range("I141").Value = 5619
$alc.Range("J141").Value = 4590
$alc.Range("K141").Value = 16857
$alc.Range("L141").Value = 13770
$alc.Range("M141").Value = -11677
$alc.Range("N141").Value = -24130

$arm = $wb.Worksheets.Item("ARM")
# Row 40 (ARM)
$arm.Range("H40").Value = 12000
$arm.Range("J40").Value = 12000
$arm.Range("L40").Value = 12000
$arm.Range("N40").Value = -12352

# Row 44 (ARM)
$arm.Range("H44").Value = 12849.857
$arm.Range("J44").Value = 12824.833
$arm.Range("L44").Value = 12824.833
$arm.Range("N44").Value = -13800.833

# Row 55 (ARM)
$arm.Range("H55").Value = 12122.857
$arm.Range("J55").Value = 12122.857
$arm.Range("L55").Value = 12122.857
$arm.Range("N55").Value = -12752.857

# Row 74 (ARM)
$arm.Range("H74").Value = 1223.0625
$arm.Range("I74").Value = 1148
$arm.Range("K74").Value = 1148
$arm.Range("M74").Value = -274
$arm.Range("N74").ClearContents()

# Row 77 (ARM)
$arm.Range("H77").Value = 1223.0625
$arm.Range("I77").Value = 1148
$arm.Range("K77").Value = 5740
$arm.Range("M77").Value = -1372
$arm.Range("N77").ClearContents()

# Row 80 (ARM)
$arm.Range("H80").Value = 25615.455
$arm.Range("J80").Value = 25615.455
$arm.Range("L80").Value = 25615.455
$arm.Range("N80").Value = -27611.455

# Row 83 (ARM)
$arm.Range("H83").Value = 25615.455
$arm.Range("J83").Value = 25615.455
$arm.Range("L83").Value = 76846.36500000001
$arm.Range("N83").Value = -86830.36500000001

$crp = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$crp.Range("H31").Value = 30415.846
$crp.Range("I31").Value = 1375.862
$crp.Range("J31").Value = 67031.48
$crp.Range("K31").Value = 1375.862
$crp.Range("L31").Value = 67031.48
$crp.Range("M31").Value = -1080.862
$crp.Range("N31").Value = -67621.48

# Row 34 (CRP)
$crp.Range("H34").Value = 30415.846
$crp.Range("I34").Value = 1375.862
$crp.Range("J34").Value = 67031.48
$crp.Range("K34").Value = 1375.862
$crp.Range("L34").Value = 67031.48
$crp.Range("M34").Value = -1173.862
$crp.Range("N34").Value = -67435.48

# Row 36 (CRP)
$crp.Range("H36").Value = 8526.5
$crp.Range("J36").Value = 15053
$crp.Range("L36").Value = 15053
$crp.Range("N36").Value = -15829

# Row 40 (CRP)
$crp.Range("H40").Value = 8526.5
$crp.Range("J40").Value = 15053
$crp.Range("L40").Value = 15053
$crp.Range("N40").Value = -15373

# Row 99 (CRP)
$crp.Range("H99").Value = 7328.2383
$crp.Range("I99").Value = 2356
$crp.Range("J99").Value = 8882.0625
$crp.Range("K99").Value = 2356
$crp.Range("L99").Value = 8882.0625
$crp.Range("M99").Value = -858
$crp.Range("N99").Value = -11878.0625

# Row 126 (CRP)
$crp.Range("H126").Value = 7328.2383
$crp.Range("I126").Value = 2356
$crp.Range("J126").Value = 8882.0625
$crp.Range("K126").Value = 7068
$crp.Range("L126").Value = 26646.1875
$crp.Range("M126").Value = -4598
$crp.Range("N126").Value = -31586.1875

# Row 132 (CRP)
$crp.Range("H132").Value = 2212.077
$crp.Range("I132").Value = 1952.8788
$crp.Range("J132").Value = 3637.6667
$crp.Range("K132").Value = 5858.636399999999
$crp.Range("L132").Value = 10913.0001
$crp.Range("M132").Value = -3328.636399999999
$crp.Range("N132").Value = -15973.0001

# Row 141 (CRP)
$crp.Range("H141").Value = 87328.28999999999
$crp.Range("J141").Value = 70236
$crp.Range("L141").Value = 70236
$crp.Range("N141").Value = -80596

$cul = $wb.Worksheets.Item("CUL")
# Row 55 (CUL)
$cul.Range("H55").Value = 10021.462
$cul.Range("I55").Value = 20376
$cul.Range("J55").Value = 3549.875
$cul.Range("K55").Value = 61128
$cul.Range("L55").Value = 10649.625
$cul.Range("M55").Value = -60951
$cul.Range("N55").Value = -11003.625

# Row 56 (CUL)
$cul.Range("H56").Value = 4385.4
$cul.Range("I56").Value = 4385.4
$cul.Range("K56").Value = 4385.4
$cul.Range("M56").Value = -3855.4

# Row 70 (CUL)
$cul.Range("H70").Value = 68960
$cul.Range("J70").Value = 2987.7778
$cul.Range("L70").Value = 8963.3334
$cul.Range("N70").Value = -9593.3334

# Row 73 (CUL)
$cul.Range("H73").Value = 68960
$cul.Range("J73").Value = 2987.7778
$cul.Range("L73").Value = 8963.3334
$cul.Range("N73").Value = -11147.3334

# Row 80 (CUL)
$cul.Range("H80").Value = 1704.5385
$cul.Range("J80").Value = 1704.5385
$cul.Range("L80").Value = 5113.6155
$cul.Range("N80").Value = -6985.6155

# Row 83 (CUL)
$cul.Range("H83").Value = 1704.5385
$cul.Range("J83").Value = 1704.5385
$cul.Range("L83").Value = 15340.8465
$cul.Range("N83").Value = -24700.8465

# Row 113 (CUL)
$cul.Range("H113").Value = 1008.5455
$cul.Range("I113").Value = 1411.8182
$cul.Range("J113").Value = 605.2727
$cul.Range("K113").Value = 4235.4546
$cul.Range("L113").Value = 1815.8181
$cul.Range("M113").Value = -2065.4546
$cul.Range("N113").Value = -6155.8181

$gsm = $wb.Worksheets.Item("GSM")
# Row 80 (GSM)
$gsm.Range("H80").Value = 143175760
$gsm.Range("I80").Value = 250551250
$gsm.Range("J80").Value = 8433.333000000001
$gsm.Range("K80").Value = 250551250
$gsm.Range("L80").Value = 8433.333000000001
$gsm.Range("M80").Value = -250550252
$gsm.Range("N80").Value = -10429.333

# Row 83 (GSM)
$gsm.Range("H83").Value = 143175760
$gsm.Range("I83").Value = 250551250
$gsm.Range("J83").Value = 8433.333000000001
$gsm.Range("K83").Value = 1252756250
$gsm.Range("L83").Value = 42166.665
$gsm.Range("M83").Value = -1252751258
$gsm.Range("N83").Value = -52150.665

# Row 122 (GSM)
$gsm.Range("H122").Value = 10000
$gsm.Range("I122").Value = 10000
$gsm.Range("J122").Value = 0
$gsm.Range("K122").Value = 30000
$gsm.Range("L122").Value = 0
$gsm.Range("M122").Value = -27550
$gsm.Range("N122").ClearContents()

$ltw = $wb.Worksheets.Item("LTW")
# Row 68 (LTW)
$ltw.Range("H68").Value = 1913.0834
$ltw.Range("J68").Value = 2578.4285
$ltw.Range("L68").Value = 2578.4285
$ltw.Range("N68").Value = -4076.4285

# Row 71 (LTW)
$ltw.Range("H71").Value = 1913.0834
$ltw.Range("J71").Value = 2578.4285
$ltw.Range("L71").Value = 12892.1425
$ltw.Range("N71").Value = -20380.1425

# Row 122 (LTW)
$ltw.Range("H122").Value = 0
$ltw.Range("I122").Value = 0
$ltw.Range("J122").Value = 0
$ltw.Range("K122").Value = 0
$ltw.Range("L122").Value = 0
$ltw.Range("M122").ClearContents()
$ltw.Range("N122").ClearContents()

# Row 132 (LTW)
$ltw.Range("H132").Value = 4666.96
$ltw.Range("I132").Value = 7421.6
$ltw.Range("J132").Value = 2830.5334
$ltw.Range("K132").Value = 22264.8
$ltw.Range("L132").Value = 8491.600199999999
$ltw.Range("M132").Value = -19734.8
$ltw.Range("N132").Value = -13551.6002

$wvr = $wb.Worksheets.Item("WVR")
# Row 54 (WVR)
$wvr.Range("H54").Value = 6382.3335
$wvr.Range("I54").Value = 6035
$wvr.Range("J54").Value = 7077
$wvr.Range("K54").Value = 6035
$wvr.Range("L54").Value = 7077
$wvr.Range("M54").Value = -5515
$wvr.Range("N54").Value = -8117

# Row 98 (WVR)
$wvr.Range("H98").Value = 0
$wvr.Range("J98").Value = 0
$wvr.Range("L98").Value = -34990
$wvr.Range("N98").ClearContents()
